$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Update the "总计" (summary) sheet: insert a new 2022-Q4 row
#    at the top of the data, push the old last row (2020-Q4) down
#    into a new row 8, and shift everything else down by one row.
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Duplicate the last existing data row (row 7 -> row 8), which both
# grows the used range to A1:D8 and carries the row style along so
# the new last row keeps the same formatting as its neighbours.
$summary.Range("A7:D7").Copy($summary.Range("A8:D8"))
$summary.Range("A8").Value = 6

# Shift the quarterly figures down by one row each (B/C/D only --
# column A is just the running 0-based index and already reads
# 0,1,2,3,4,5 top to bottom, so it does not need to change).
$summary.Range("B2").Value = '2022-Q4'
$summary.Range("C2").Value = 39
$summary.Range("D2").Value = 12.58
$summary.Range("B3").Value = '2022-Q3'
$summary.Range("C3").Value = 22
$summary.Range("D3").Value = 10.3
$summary.Range("B4").Value = '2022-Q2'
$summary.Range("C4").Value = 18
$summary.Range("D4").Value = 14.49
$summary.Range("B5").Value = '2022-Q1'
$summary.Range("C5").Value = 12
$summary.Range("D5").Value = 18.3
$summary.Range("B6").Value = '2021-Q4'
$summary.Range("C6").Value = 16
$summary.Range("D6").Value = 11.34
$summary.Range("B7").Value = '2021-Q3'
$summary.Range("C7").Value = 10
$summary.Range("D7").Value = 9.66

# ---------------------------------------------------------------
# 2) Insert a brand-new "2022-Q4" sheet right after "总计" (i.e. as
#    the new 2nd sheet), holding the per-fund holdings detail for
#    that quarter. Every other sheet keeps its name and just
#    shifts one position to the right/down.
# ---------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $summary)
$q4.Name = "2022-Q4"

# Header row (bold, centered, top-aligned, boxed -- matches the
# header styling used on every other quarterly sheet).
$headerRange = $q4.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$q4.Range("B1").Value = '基金代码'
$q4.Range("C1").Value = '基金名称'
$q4.Range("D1").Value = '基金规模'
$q4.Range("E1").Value = '股票总仓位'
$q4.Range("F1").Value = '仓位占比'
$q4.Range("G1").Value = '持有市值(亿元)'
$q4.Range("H1").Value = '仓位排名'

# Column A (the running index) is bold/centered/top/boxed on every
# other sheet too, so format the whole A2:A40 block up front.
$indexRange = $q4.Range("A2:A40")
$indexRange.Font.Bold = $true
$indexRange.HorizontalAlignment = -4108
$indexRange.VerticalAlignment = -4160
$indexRange.Borders.LineStyle = 1

# Columns B-G hold text (fund codes/names and percentage-like
# figures that must keep leading zeros / fixed decimals as text),
# so force the Text number format before writing the strings.
$q4.Range("B2:G40").NumberFormat = "@"

# Per-fund holdings rows.
$q4.Range("A2").Value = 0
$q4.Range("B2").Value = '011531'
$q4.Range("C2").Value = '朱雀恒心一年持有期混合'
$q4.Range("D2").Value = '58.59'
$q4.Range("E2").Value = '93.78'
$q4.Range("F2").Value = '6.05'
$q4.Range("G2").Value = '3.5447'
$q4.Range("H2").Value = 4
$q4.Range("A3").Value = 1
$q4.Range("B3").Value = '010141'
$q4.Range("C3").Value = '朱雀企业优选股票A'
$q4.Range("D3").Value = '27.66'
$q4.Range("E3").Value = '93.07'
$q4.Range("F3").Value = '6.76'
$q4.Range("G3").Value = '1.8698'
$q4.Range("H3").Value = 3
$q4.Range("A4").Value = 2
$q4.Range("B4").Value = '007493'
$q4.Range("C4").Value = '朱雀产业臻选混合A'
$q4.Range("D4").Value = '34.98'
$q4.Range("E4").Value = '93.77'
$q4.Range("F4").Value = '4.69'
$q4.Range("G4").Value = '1.6406'
$q4.Range("H4").Value = 4
$q4.Range("A5").Value = 3
$q4.Range("B5").Value = '010922'
$q4.Range("C5").Value = '朱雀匠心一年持有期混合'
$q4.Range("D5").Value = '15.00'
$q4.Range("E5").Value = '93.80'
$q4.Range("F5").Value = '7.05'
$q4.Range("G5").Value = '1.0575'
$q4.Range("H5").Value = 3
$q4.Range("A6").Value = 4
$q4.Range("B6").Value = '050009'
$q4.Range("C6").Value = '博时新兴成长混合'
$q4.Range("D6").Value = '23.52'
$q4.Range("E6").Value = '90.03'
$q4.Range("F6").Value = '4.44'
$q4.Range("G6").Value = '1.0443'
$q4.Range("H6").Value = 6
$q4.Range("A7").Value = 5
$q4.Range("B7").Value = '007494'
$q4.Range("C7").Value = '朱雀产业臻选混合C'
$q4.Range("D7").Value = '12.43'
$q4.Range("E7").Value = '93.77'
$q4.Range("F7").Value = '4.69'
$q4.Range("G7").Value = '0.5830'
$q4.Range("H7").Value = 4
$q4.Range("A8").Value = 6
$q4.Range("B8").Value = '008294'
$q4.Range("C8").Value = '朱雀企业优胜股票A'
$q4.Range("D8").Value = '11.44'
$q4.Range("E8").Value = '93.93'
$q4.Range("F8").Value = '4.68'
$q4.Range("G8").Value = '0.5354'
$q4.Range("H8").Value = 4
$q4.Range("A9").Value = 7
$q4.Range("B9").Value = '010142'
$q4.Range("C9").Value = '朱雀企业优选股票C'
$q4.Range("D9").Value = '5.26'
$q4.Range("E9").Value = '93.07'
$q4.Range("F9").Value = '6.76'
$q4.Range("G9").Value = '0.3556'
$q4.Range("H9").Value = 3
$q4.Range("A10").Value = 8
$q4.Range("B10").Value = '910005'
$q4.Range("C10").Value = '东方红启兴三年持有混合'
$q4.Range("D10").Value = '5.75'
$q4.Range("E10").Value = '92.69'
$q4.Range("F10").Value = '5.00'
$q4.Range("G10").Value = '0.2875'
$q4.Range("H10").Value = 9
$q4.Range("A11").Value = 9
$q4.Range("B11").Value = '007040'
$q4.Range("C11").Value = '新疆前海联合泳隆灵活配置混合C'
$q4.Range("D11").Value = '7.66'
$q4.Range("E11").Value = '93.78'
$q4.Range("F11").Value = '3.54'
$q4.Range("G11").Value = '0.2712'
$q4.Range("H11").Value = 9
$q4.Range("A12").Value = 10
$q4.Range("B12").Value = '050010'
$q4.Range("C12").Value = '博时特许价值混合'
$q4.Range("D12").Value = '5.60'
$q4.Range("E12").Value = '91.87'
$q4.Range("F12").Value = '4.76'
$q4.Range("G12").Value = '0.2666'
$q4.Range("H12").Value = 4
$q4.Range("A13").Value = 11
$q4.Range("B13").Value = '015729'
$q4.Range("C13").Value = '朱雀碳中和三年持有期混合'
$q4.Range("D13").Value = '3.60'
$q4.Range("E13").Value = '49.87'
$q4.Range("F13").Value = '5.61'
$q4.Range("G13").Value = '0.2020'
$q4.Range("H13").Value = 2
$q4.Range("A14").Value = 12
$q4.Range("B14").Value = '000534'
$q4.Range("C14").Value = '长盛高端装备制造灵活配置混合A'
$q4.Range("D14").Value = '6.28'
$q4.Range("E14").Value = '88.01'
$q4.Range("F14").Value = '2.85'
$q4.Range("G14").Value = '0.1790'
$q4.Range("H14").Value = 10
$q4.Range("A15").Value = 13
$q4.Range("B15").Value = '007880'
$q4.Range("C15").Value = '朱雀产业智选混合A'
$q4.Range("D15").Value = '4.14'
$q4.Range("E15").Value = '92.62'
$q4.Range("F15").Value = '4.11'
$q4.Range("G15").Value = '0.1702'
$q4.Range("H15").Value = 5
$q4.Range("A16").Value = 14
$q4.Range("B16").Value = '008295'
$q4.Range("C16").Value = '朱雀企业优胜股票C'
$q4.Range("D16").Value = '2.44'
$q4.Range("E16").Value = '93.93'
$q4.Range("F16").Value = '4.68'
$q4.Range("G16").Value = '0.1142'
$q4.Range("H16").Value = 4
$q4.Range("A17").Value = 15
$q4.Range("B17").Value = '001215'
$q4.Range("C17").Value = '博时沪港深优质企业混合A'
$q4.Range("D17").Value = '2.08'
$q4.Range("E17").Value = '91.73'
$q4.Range("F17").Value = '4.52'
$q4.Range("G17").Value = '0.0940'
$q4.Range("H17").Value = 5
$q4.Range("A18").Value = 16
$q4.Range("B18").Value = '005903'
$q4.Range("C18").Value = '泰达宏利绩优增长灵活配置混合A'
$q4.Range("D18").Value = '1.85'
$q4.Range("E18").Value = '87.04'
$q4.Range("F18").Value = '3.32'
$q4.Range("G18").Value = '0.0614'
$q4.Range("H18").Value = 8
$q4.Range("A19").Value = 17
$q4.Range("B19").Value = '015576'
$q4.Range("C19").Value = '泰达宏利绩优增长灵活配置混合C'
$q4.Range("D19").Value = '1.70'
$q4.Range("E19").Value = '87.04'
$q4.Range("F19").Value = '3.32'
$q4.Range("G19").Value = '0.0564'
$q4.Range("H19").Value = 8
$q4.Range("A20").Value = 18
$q4.Range("B20").Value = '000598'
$q4.Range("C20").Value = '长盛生态环境主题灵活配置混合'
$q4.Range("D20").Value = '1.88'
$q4.Range("E20").Value = '85.74'
$q4.Range("F20").Value = '2.83'
$q4.Range("G20").Value = '0.0532'
$q4.Range("H20").Value = 9
$q4.Range("A21").Value = 19
$q4.Range("B21").Value = '001892'
$q4.Range("C21").Value = '长盛新兴成长主题灵活配置混合'
$q4.Range("D21").Value = '1.14'
$q4.Range("E21").Value = '88.71'
$q4.Range("F21").Value = '3.66'
$q4.Range("G21").Value = '0.0417'
$q4.Range("H21").Value = 7
$q4.Range("A22").Value = 20
$q4.Range("B22").Value = '004128'
$q4.Range("C22").Value = '新疆前海联合泳隆灵活配置混合A'
$q4.Range("D22").Value = '0.84'
$q4.Range("E22").Value = '93.78'
$q4.Range("F22").Value = '3.54'
$q4.Range("G22").Value = '0.0297'
$q4.Range("H22").Value = 9
$q4.Range("A23").Value = 21
$q4.Range("B23").Value = '007881'
$q4.Range("C23").Value = '朱雀产业智选混合C'
$q4.Range("D23").Value = '0.71'
$q4.Range("E23").Value = '92.62'
$q4.Range("F23").Value = '4.11'
$q4.Range("G23").Value = '0.0292'
$q4.Range("H23").Value = 5
$q4.Range("A24").Value = 22
$q4.Range("B24").Value = '080002'
$q4.Range("C24").Value = '长盛创新先锋混合A'
$q4.Range("D24").Value = '0.73'
$q4.Range("E24").Value = '77.72'
$q4.Range("F24").Value = '2.81'
$q4.Range("G24").Value = '0.0205'
$q4.Range("H24").Value = 6
$q4.Range("A25").Value = 23
$q4.Range("B25").Value = '001261'
$q4.Range("C25").Value = '中融新机遇灵活配置混合'
$q4.Range("D25").Value = '0.46'
$q4.Range("E25").Value = '78.82'
$q4.Range("F25").Value = '3.67'
$q4.Range("G25").Value = '0.0169'
$q4.Range("H25").Value = 9
$q4.Range("A26").Value = 24
$q4.Range("B26").Value = '002156'
$q4.Range("C26").Value = '长盛盛世灵活配置混合A'
$q4.Range("D26").Value = '1.66'
$q4.Range("E26").Value = '27.33'
$q4.Range("F26").Value = '0.87'
$q4.Range("G26").Value = '0.0144'
$q4.Range("H26").Value = 9
$q4.Range("A27").Value = 25
$q4.Range("B27").Value = '003235'
$q4.Range("C27").Value = '信诚至利灵活配置混合C'
$q4.Range("D27").Value = '1.80'
$q4.Range("E27").Value = '23.20'
$q4.Range("F27").Value = '0.50'
$q4.Range("G27").Value = '0.0090'
$q4.Range("H27").Value = 8
$q4.Range("A28").Value = 26
$q4.Range("B28").Value = '001402'
$q4.Range("C28").Value = '信诚新选回报灵活配置混合A'
$q4.Range("D28").Value = '1.06'
$q4.Range("E28").Value = '24.28'
$q4.Range("F28").Value = '0.52'
$q4.Range("G28").Value = '0.0055'
$q4.Range("H28").Value = 8
$q4.Range("A29").Value = 27
$q4.Range("B29").Value = '970046'
$q4.Range("C29").Value = '东海证券海睿健行灵活配置混合A'
$q4.Range("D29").Value = '0.14'
$q4.Range("E29").Value = '87.16'
$q4.Range("F29").Value = '3.41'
$q4.Range("G29").Value = '0.0048'
$q4.Range("H29").Value = 9
$q4.Range("A30").Value = 28
$q4.Range("B30").Value = '002555'
$q4.Range("C30").Value = '博时沪港深优质企业混合C'
$q4.Range("D30").Value = '0.08'
$q4.Range("E30").Value = '91.73'
$q4.Range("F30").Value = '4.52'
$q4.Range("G30").Value = '0.0036'
$q4.Range("H30").Value = 5
$q4.Range("A31").Value = 29
$q4.Range("B31").Value = '003234'
$q4.Range("C31").Value = '信诚至利灵活配置混合A'
$q4.Range("D31").Value = '0.65'
$q4.Range("E31").Value = '23.20'
$q4.Range("F31").Value = '0.50'
$q4.Range("G31").Value = '0.0032'
$q4.Range("H31").Value = 8
$q4.Range("A32").Value = 30
$q4.Range("B32").Value = '970047'
$q4.Range("C32").Value = '东海证券海睿健行灵活配置混合B'
$q4.Range("D32").Value = '0.09'
$q4.Range("E32").Value = '87.16'
$q4.Range("F32").Value = '3.41'
$q4.Range("G32").Value = '0.0031'
$q4.Range("H32").Value = 9
$q4.Range("A33").Value = 31
$q4.Range("B33").Value = '004157'
$q4.Range("C33").Value = '信诚至诚灵活配置混合A'
$q4.Range("D33").Value = '0.56'
$q4.Range("E33").Value = '23.32'
$q4.Range("F33").Value = '0.50'
$q4.Range("G33").Value = '0.0028'
$q4.Range("H33").Value = 7
$q4.Range("A34").Value = 32
$q4.Range("B34").Value = '005538'
$q4.Range("C34").Value = '中航新起航灵活配置混合C'
$q4.Range("D34").Value = '0.04'
$q4.Range("E34").Value = '69.18'
$q4.Range("F34").Value = '5.69'
$q4.Range("G34").Value = '0.0023'
$q4.Range("H34").Value = 8
$q4.Range("A35").Value = 33
$q4.Range("B35").Value = '002030'
$q4.Range("C35").Value = '信诚新选回报灵活配置混合B'
$q4.Range("D35").Value = '0.36'
$q4.Range("E35").Value = '24.28'
$q4.Range("F35").Value = '0.52'
$q4.Range("G35").Value = '0.0019'
$q4.Range("H35").Value = 8
$q4.Range("A36").Value = 34
$q4.Range("B36").Value = '005537'
$q4.Range("C36").Value = '中航新起航灵活配置混合A'
$q4.Range("D36").Value = '0.02'
$q4.Range("E36").Value = '69.18'
$q4.Range("F36").Value = '5.69'
$q4.Range("G36").Value = '0.0011'
$q4.Range("H36").Value = 8
$q4.Range("A37").Value = 35
$q4.Range("B37").Value = '004158'
$q4.Range("C37").Value = '信诚至诚灵活配置混合B'
$q4.Range("D37").Value = '0.16'
$q4.Range("E37").Value = '23.32'
$q4.Range("F37").Value = '0.50'
$q4.Range("G37").Value = '0.0008'
$q4.Range("H37").Value = 7
$q4.Range("A38").Value = 36
$q4.Range("B38").Value = '002157'
$q4.Range("C38").Value = '长盛盛世灵活配置混合C'
$q4.Range("D38").Value = '0.04'
$q4.Range("E38").Value = '27.33'
$q4.Range("F38").Value = '0.87'
$q4.Range("G38").Value = '0.0003'
$q4.Range("H38").Value = 9
$q4.Range("A39").Value = 37
$q4.Range("B39").Value = '012716'
$q4.Range("C39").Value = '长盛创新先锋混合C'
$q4.Range("D39").Value = '0.00'
$q4.Range("E39").Value = '77.72'
$q4.Range("F39").Value = '2.81'
$q4.Range("G39").NumberFormat = "General"
$q4.Range("G39").Value = 0
$q4.Range("H39").Value = 6
$q4.Range("A40").Value = 38
$q4.Range("B40").Value = '017485'
$q4.Range("C40").Value = '长盛高端装备制造灵活配置混合C'
$q4.Range("D40").Value = '0.00'
$q4.Range("E40").Value = '88.01'
$q4.Range("F40").Value = '2.85'
$q4.Range("G40").NumberFormat = "General"
$q4.Range("G40").Value = 0
$q4.Range("H40").Value = 10
